# Testing 2.xlsx - "Updated Testing2 document after testing the test cases"
#
# The test cases in Sheet1 (rows 4-9) get a "Tester 1" name filled in
# (column D) and a "Pass" result recorded (column E, highlighted green
# like the existing "Result"/"Tester 2" columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Green fill used by the existing "Pass" cells (column H) - RGB(146,208,80)
$passFillColor = 5296274

foreach ($r in 4..9) {
    $ws.Range("D$r").Value = "Aneesh Dalvi"

    $ws.Range("E$r").Value = "Pass"
    $ws.Range("E$r").Interior.Color = $passFillColor
}

# Row heights grew now that columns D & E are populated (re-autofit).
$ws.Rows.Item(1).RowHeight = 93.75
$ws.Rows.Item(3).RowHeight = 30
$ws.Rows.Item(4).RowHeight = 45
$ws.Rows.Item(5).RowHeight = 90
$ws.Rows.Item(6).RowHeight = 75
$ws.Rows.Item(7).RowHeight = 75
$ws.Rows.Item(8).RowHeight = 45
$ws.Rows.Item(9).RowHeight = 105

# Selection left on the newly-filled-in Result column after testing.
$ws.Range("E4:E9").Select() | Out-Null
